$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add two new header labels next to the existing ones ---
$ws.Range("C1").Value = "table_id"
$ws.Range("D1").Value = "table_class"

# --- Row 2: replace the old "http://localhost" sample with a real local URL ---
# Drop the old (rich-text styled) hyperlink first so we don't end up with stale formatting
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Value = "http://127.0.0.1:8000"
$ws.Hyperlinks.Add($ws.Range("A2"), "http://127.0.0.1:8000")
$ws.Range("B2").Value = "firefox"

# --- Row 3: a brand-new sample site, including the scraping-specific columns ---
$ws.Range("A3").Value = "https://www.vgchartz.com/charts/platform_totals/Hardware.php"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.vgchartz.com/charts/platform_totals/Hardware.php")
$ws.Range("B3").Value = "firefox"
$ws.Range("C3").Value = "myTable"
$ws.Range("D3").Value = "chart"

# Carry the hyperlink-style formatting down into the still-empty rows below,
# matching the look of the populated URL cells above
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4:B5").PasteSpecial(-4122)

# --- Column A needs to be noticeably wider to fit the longer URLs ---
$ws.Columns("A").ColumnWidth = 31.6

# --- Leave the selection where the author last left it ---
$ws.Range("F7").Select() | Out-Null
